$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-11-24"

# Update the row-13 label ("November (through 11-23)" -> "11-24")
$ws.Range("A13").Value = "November (through 11-24)"

# --- Row 13 updates ---
$ws.Range("C13").Value = 23
$ws.Range("D13").Value = 0.0417

$ws.Range("F13").Value = 53
$ws.Range("G13").Value = 0.1167

$ws.Range("I13").Value = 89
$ws.Range("J13").Value = 0.022

$ws.Range("K13").Value = 7
$ws.Range("L13").Value = 41
$ws.Range("M13").Value = 0.1458

$ws.Range("O13").Value = 37
$ws.Range("P13").Value = 0.1395

$ws.Range("R13").Value = 166
$ws.Range("S13").Value = 0.046

$ws.Range("U13").Value = 162
$ws.Range("V13").Value = 0.0182

# --- Row 14 (Total) updates ---
$ws.Range("C14").Value = 249
$ws.Range("D14").Value = 0.117

$ws.Range("F14").Value = 487
$ws.Range("G14").Value = 0.1081

$ws.Range("I14").Value = 738
$ws.Range("J14").Value = 0.07870000000000001

$ws.Range("K14").Value = 73
$ws.Range("L14").Value = 590
$ws.Range("M14").Value = 0.1101

$ws.Range("O14").Value = 471
$ws.Range("P14").Value = 0.1029

$ws.Range("R14").Value = 1169
$ws.Range("S14").Value = 0.0504

$ws.Range("U14").Value = 1513
$ws.Range("V14").Value = 0.0591
